$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3058

$ws.Range("A3").Value = "granodiorite"
$ws.Range("B3").Value = 908

$ws.Range("A4").Value = "quartz monzodiorite`nquartz monzogabbro"
$ws.Range("B4").Value = 394

$ws.Range("A5").Value = "quartz monzonite"
$ws.Range("B5").Value = 199

$ws.Range("A6").Value = "monzodiorite monzogabbro"
$ws.Range("B6").Value = 30

$ws.Range("A7").Value = "syeno granite"
$ws.Range("B7").Value = 27

$ws.Range("A8").Value = "tonalite"
$ws.Range("B8").Value = 26

$ws.Range("A10").Value = "quartz diorite`nquartz gabbro`nquartz anorthosite"
$ws.Range("B10").Value = 5

$ws.Range("A11").Value = "monzonite"
$ws.Range("B11").Value = 5

$ws.Range("A12").Value = "diorite gabbro anorthosite"
$ws.Range("B12").Value = 2
